# Applies the "update database and change read_price algorithm" commit:
#  - company name (B5) gains a "ص." abbreviation
#  - the oldest reporting period (1396/12) is dropped and a new one
#    (1401/12) is appended, so the five year-header columns (E:I) now
#    read 1397/12 .. 1401/12 instead of 1396/12 .. 1400/12
#  - every data column shifts one year to the left (old F->E, G->F,
#    H->G, I->H) and a freshly reported value is filled into column I

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Company name -----------------------------------------------------
$ws.Range("B5").Value = "کیمیا-ص. معدنی کیمیای زنجان گستران"

# --- Year header row labels (used for both header rows, 8 and 24) -----
$ws.Range("E8").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E24").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F24").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G24").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H24").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I24").Value = "دوازده ماهه منتهی به 1401/12"

# --- Row 10: هزینه حمل و نقل و انتقال ----------------------------------
$ws.Range("E10").Value = 3831
$ws.Range("F10").Value = 631
$ws.Range("G10").Value = 0
# H10, I10 unchanged (0, 0)

# --- Row 13: حق العمل و کمیسیون فروش -----------------------------------
$ws.Range("E13").Value = 1787
$ws.Range("F13").Value = 659
$ws.Range("G13").Value = 2564
$ws.Range("H13").Value = 6820
$ws.Range("I13").Value = 1350

# --- Row 14: هزینه تبلیغات ---------------------------------------------
$ws.Range("E14").Value = 357
$ws.Range("F14").Value = 0
# G14, H14, I14 unchanged (0, 0, 0)

# --- Row 15: هزینه مواد مصرفی -------------------------------------------
$ws.Range("E15").Value = 66
$ws.Range("F15").Value = 315
$ws.Range("G15").Value = 695
$ws.Range("H15").Value = 1513
$ws.Range("I15").Value = 1599

# --- Row 16: هزینه انرژی (آب، برق، گاز و سوخت) ---------------------------
$ws.Range("E16").Value = 446
$ws.Range("F16").Value = 625
$ws.Range("G16").Value = 810
$ws.Range("H16").Value = 3491
$ws.Range("I16").Value = 6975

# --- Row 17: هزینه استهلاک ----------------------------------------------
$ws.Range("E17").Value = 15145
$ws.Range("F17").Value = 19580
$ws.Range("G17").Value = 29131
$ws.Range("H17").Value = 76076
$ws.Range("I17").Value = 97205

# --- Row 19: سایر هزینه ها ----------------------------------------------
$ws.Range("E19").Value = 7934
$ws.Range("F19").Value = 12131
$ws.Range("G19").Value = 40658
$ws.Range("H19").Value = 229117
$ws.Range("I19").Value = 255308

# --- Row 20: جمع (totals row) -------------------------------------------
$ws.Range("E20").Value = 29566
$ws.Range("F20").Value = 33941
$ws.Range("G20").Value = 73858
$ws.Range("H20").Value = 317017
$ws.Range("I20").Value = 362437

# --- Row 26: تعداد پرسنل غیر تولیدی شرکت --------------------------------
$ws.Range("E26").Value = 113
$ws.Range("F26").Value = 110
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 27
# I26 unchanged (27)

# --- Row 27: تعداد پرسنل تولیدی شرکت ------------------------------------
$ws.Range("E27").Value = 71
$ws.Range("F27").Value = 70
$ws.Range("G27").Value = 141
$ws.Range("H27").Value = 150
$ws.Range("I27").Value = 156
